$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "descr"
$ws.Range("E1").Value = "is_active"

# Data rows: lang_code, code, name, descr, is_active
$data = @(
    @("eng", "UIN",      "Unique Identification Number",   "National ID given to the applicant", $true),
    @("eng", "PRID",     "Pre-registration ID",             "ID assigned after Pre-registration", $true),
    @("eng", "RID",      "Registration ID",                 "ID assigned after registration", $true),
    @("eng", "VID",      "Virtual ID",                      "ID used in replacement of UIN", $true),
    @("eng", "Token ID", "Token ID",                        "ID used by a vendor for an applicant", $true),
    @("fra", "UIN",      "Numéro didentification unique",   "Carte didentité nationale fournie au demandeur", $true),
    @("fra", "PRID",     "ID de pré-inscription",           "ID attribué après la pré-inscription", $true),
    @("fra", "RID",      "ID denregistrement",              "ID attribué après lenregistrement", $true),
    @("fra", "VID",      "ID virtuel",                      "Identifiant utilisé en remplacement de UIN", $true),
    @("fra", "Token ID", "ID de jeton",                     "ID utilisé par un fournisseur pour un demandeur", $true)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Column A (lang_code) in data rows uses the same style as the header row (s="1")
$ws.Range("A1").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
